$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# labelName (column C) updates first
$ws.Range("C11").Value = "Eggs Sold (Y/N)"
$ws.Range("C14").Value = "Eggs Sold (Qty)"

# prettyName (column D) updates next
$ws.Range("D11").Value = "Household Sold Eggs"
$ws.Range("D14").Value = "Number of Chicken Eggs Sold"

# Update the active cell selection to match the saved workbook state
$ws.Range("E26").Select()
